$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = -0.0062038245579894317
$ws.Range("B1").Value = -0.0061970331349624477
$ws.Range("A2").Value = -0.010780502056723006
$ws.Range("B2").Value = -0.010780285198182735
$ws.Range("A3").Value = -0.01706229068708795
$ws.Range("B3").Value = -0.017066299336157158
$ws.Range("A4").Value = -0.073732523962045979
$ws.Range("B4").Value = -0.073730918118873376
